# Updated unit test plan
# Fill in the previously-blank "Preconditions" / "Method Inputs" / "Expected Result"
# columns (E:G) for the __init__ test cases (rows 7-10) and the attribute-getter
# test cases (rows 11-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$none = "None"
$attrsSet = "Attributes are set"
$valueError = "Valuerror"
$inputs = "Title = `"Atomic Habits`"`nauthor = `"James Clear`"`nGenre = `"NON_FICTION`""

# Row 7 - first writes establish the shared-string order: None, Attributes are set,
# Valuerror, then the multi-line Method Inputs text.
$ws.Range("E7").Value = $none
$ws.Range("G7").Value = $attrsSet
$ws.Range("G8").Value = $valueError
$ws.Range("F7").Value = $inputs

# Row 8 remainder
$ws.Range("F8").Value = $inputs
$ws.Range("E8").Value = $none

# Row 9
$ws.Range("G9").Value = $valueError
$ws.Range("F9").Value = $inputs
$ws.Range("E9").Value = $none

# Row 10
$ws.Range("G10").Value = $valueError
$ws.Range("F10").Value = $inputs
$ws.Range("E10").Value = $none

# Rows 11-13: attribute getters - Preconditions column (E) gets the input text,
# Method Inputs column (F) gets "None"; Expected Result (G) stays blank.
$ws.Range("E11").Value = $inputs
$ws.Range("F11").Value = $none

$ws.Range("E12").Value = $inputs
$ws.Range("F12").Value = $none

$ws.Range("E13").Value = $inputs
$ws.Range("F13").Value = $none

# Update the saved selection/view to match the author's last position.
$ws.Range("G11").Select() | Out-Null
